$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update D29 with the new dask exercise link text.
# Use Formula with a leading apostrophe so the cell keeps its existing
# "quote prefix" text style (matches how the original author typed
# leading-dash text into this column).
$ws.Range("D29").Formula = "'" + '- `Link <exercises/Exercise_dask_realdata.ipynb>`_'

# Add new B30 cell text "Catchup Time"
$ws.Range("B30").Value = "Catchup Time"

# Update the active selection to B31 (matches the authored diff's UI state)
$ws.Range("B31").Select()
